$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.836.92"
$ws.Range("E2").Value = "'  -0.23%  "
$ws.Range("D3").Value = "'1.636.37"
$ws.Range("E3").Value = "'  +0.11%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("E5").Value = "'  +0.67%  "
$ws.Range("D6").Value = "'0.5058"
$ws.Range("E6").Value = "'  +0.08%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "'  +0.16%  "
$ws.Range("E8").Value = "'  +0.10%  "
$ws.Range("D9").Value = "'0.06425"
$ws.Range("E9").Value = "'  +1.15%  "
$ws.Range("D10").Value = "'19.54"
$ws.Range("E10").Value = "'  -0.71%  "
$ws.Range("D11").Value = "'0.07770"
$ws.Range("E11").Value = "'  +0.29%  "
$ws.Range("E12").Value = "'  -0.16%  "
$ws.Range("D13").Value = "'1.862.18"
$ws.Range("E13").Value = "'  +0.07%  "
$ws.Range("D14").Value = "'1.633.62"
$ws.Range("E14").Value = "'  -0.07%  "
$ws.Range("D15").Value = "'0.5622"
$ws.Range("E15").Value = "'  +3.36%  "
$ws.Range("D16").Value = "'0.0₅7591"
$ws.Range("E16").Value = "'  -1.76%  "
$ws.Range("D17").Value = "'63.10"
$ws.Range("E17").Value = "'  -1.45%  "
$ws.Range("D18").Value = "'25.841.42"
$ws.Range("E18").Value = "'  -0.29%  "
$ws.Range("E19").Value = "'  +0.11%  "
$ws.Range("D20").Value = "'194.95"
$ws.Range("E20").Value = "'  -0.19%  "
$ws.Range("E21").Value = "'  -2.53%  "
$ws.Range("D22").Value = "'9.873"
$ws.Range("E22").Value = "'  -0.38%  "
$ws.Range("E23").Value = "'  -0.33%  "
$ws.Range("E24").Value = "'  +0.01%  "
$ws.Range("D25").Value = "'1.799"
$ws.Range("E25").Value = "'  -4.76%  "
$ws.Range("B26").Value = "'Monero"
$ws.Range("C26").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'139.84"
$ws.Range("E26").Value = "'  -2.15%  "
$ws.Range("B27").Value = "'Stellar"
$ws.Range("C27").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.1267"
$ws.Range("E27").Value = "'  +2.07%  "
$ws.Range("D28").Value = "'6.774"
$ws.Range("E28").Value = "'  -0.72%  "
$ws.Range("E29").Value = "'  -1.03%  "
$ws.Range("D30").Value = "'1.243"
$ws.Range("E30").Value = "'  +0.57%  "
$ws.Range("D31").Value = "'0.04863"
$ws.Range("E31").Value = "'  -0.07%  "
$ws.Range("D32").Value = "'3.294"
$ws.Range("E32").Value = "'  +1.80%  "
$ws.Range("E33").Value = "'  +0.64%  "
$ws.Range("D34").Value = "'1.557"
$ws.Range("E34").Value = "'  +0.68%  "
$ws.Range("E35").Value = "'  +0.14%  "
$ws.Range("D36").Value = "'0.9032"
$ws.Range("E36").Value = "'  -0.83%  "
$ws.Range("D37").Value = "'2.578"
$ws.Range("E37").Value = "'  +0.28%  "
$ws.Range("D38").Value = "'1.130.92"
$ws.Range("E38").Value = "'  +0.74%  "
$ws.Range("D39").Value = "'0.5501"
$ws.Range("E39").Value = "'  -0.01%  "
$ws.Range("E40").Value = "'  +0.08%  "
$ws.Range("D41").Value = "'0.9999"
$ws.Range("E41").Value = "'  -0.13%  "
$ws.Range("D42").Value = "'5.518"
$ws.Range("E42").Value = "'  -1.32%  "
$ws.Range("D43").Value = "'0.7999"
$ws.Range("E43").Value = "'  -0.49%  "
$ws.Range("D44").Value = "'97.63"
$ws.Range("E44").Value = "'  -0.90%  "
$ws.Range("D45").Value = "'1.772.55"
$ws.Range("E45").Value = "'  +0.13%  "
$ws.Range("D46").Value = "'0.0₈113"
$ws.Range("E46").Value = "'  -8.06%  "
$ws.Range("E47").Value = "'  +0.74%  "
$ws.Range("E48").Value = "'  -2.03%  "
$ws.Range("E49").Value = "'  -2.36%  "
$ws.Range("D50").Value = "'7.674"
$ws.Range("E50").Value = "'  +2.40%  "
$ws.Range("E51").Value = "'  +0.10%  "
